# Updated symbol list on Tue Dec 20 15:09:38 UTC 2022 with GitHub Actions
#
# All data cells on the sheet are stored as text (inlineStr) in the source
# workbook, even though many of them look numeric ("248.61", "14", ...).
# Excel's COM `.Value` setter auto-converts a numeric-looking string into a
# real number, which would corrupt both the value (dropping trailing
# zeroes, e.g. "250.00" -> 250) and the cell's text type. Forcing the
# cell's NumberFormat to "@" (Text) before assigning the value keeps the
# assignment as a literal string, matching the target edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D ("Price") updates, keyed by row number.
$priceUpdates = @{
    2  = "250.00"
    3  = "22.83"
    4  = "5.413"
    5  = "0.05653"
    6  = "3.429"
    7  = "6.368"
    8  = "0.8199"
    9  = "0.9239"
    10 = "0.1433"
    11 = "0.07533"
    12 = "0.03165"
    13 = "0.03086"
    14 = "0.09316"
    15 = "3.557"
    17 = "0.04703"
    18 = "0.0005774"
    19 = "0.006372"
    20 = "0.005020"
    21 = "0.001030"
    22 = "0.0001500"
    23 = "3.725"
    24 = "2.163"
    26 = "0.1301"
    28 = "0.0002991"
    40 = "0.04002"
    41 = "0.006907"
    42 = "0.1071"
    43 = "0.002793"
    44 = "0.007565"
    45 = "0.00005556"
    46 = "0.00000000749"
    48 = "0.5982"
    49 = "0.2236"
    50 = "0.00002098"
    51 = "0.01009"
}

foreach ($row in $priceUpdates.Keys) {
    $cell = $ws.Range("D$row")
    $cell.NumberFormat = "@"
    $cell.Value = $priceUpdates[$row]
}

# Column E ("Volume(1h)") text tweaks - a couple of rows gained/lost a
# "Bestin24h"/"Worstin24h" suffix between runs.
$cellE18 = $ws.Range("E18")
$cellE18.NumberFormat = "@"
$cellE18.Value = "17OneONEWorstin24h"

$cellE27 = $ws.Range("E27")
$cellE27.NumberFormat = "@"
$cellE27.Value = "26AAXTokenAAB"

# Column G ("Hora") - the scrape hour moved from 14 to 15 for every data
# row (2 through 51).
for ($row = 2; $row -le 51; $row++) {
    $cellG = $ws.Range("G$row")
    $cellG.NumberFormat = "@"
    $cellG.Value = "15"
}
